# Add the new "2015_monthly" worksheet (monthly resident / non-resident
# termination counts) to the end of the workbook, matching the
# "add monthly data to reports" commit.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab strip (sheetId 11 / rId11 / sheet11.xml, as in the target).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2015_monthly"

# Header row.
$ws.Cells.Item(1, 1).Value = "Month"
$ws.Cells.Item(1, 2).Value = "Non-Resident Terminations Month Resident Terminations"

# Monthly data: month name, non-resident terminations, resident terminations.
$months = @("January", "February", "March", "April", "May", "June", `
            "July", "August", "September", "October", "November", "December")
$nonResident = @(753, 703, 670, 657, 681, 644, 621, 527, 581, 592, 504, 554)
$resident    = @(59, 47, 27, 44, 43, 39, 48, 23, 44, 33, 30, 33)

for ($i = 0; $i -lt $months.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $months[$i]
    $ws.Cells.Item($row, 2).Value = $nonResident[$i]
    $ws.Cells.Item($row, 3).Value = $resident[$i]
}

# Match the saved selection/active-cell state recorded for this sheet.
$ws.Range("E6").Select()
